$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 6.14
$ws.Range("I2").Value = 6.14
$ws.Range("K2").Value = 2.71
$ws.Range("M2").Value = 4.52
$ws.Range("E3").Value = 6.35
$ws.Range("H3").Value = 6.6
$ws.Range("I3").Value = 6.35
$ws.Range("J3").Value = 5.48
$ws.Range("K3").Value = 2.28
$ws.Range("L3").Value = 6.6
$ws.Range("M3").Value = 4.28
$ws.Range("E4").Value = 6.51
$ws.Range("H4").Value = 6.4
$ws.Range("I4").Value = 6.51
$ws.Range("J4").Value = 5.31
$ws.Range("K4").Value = 2.89
$ws.Range("L4").Value = 6.4
$ws.Range("M4").Value = 4.63
$ws.Range("E5").Value = 6.73
$ws.Range("H5").Value = 6.12
$ws.Range("I5").Value = 6.73
$ws.Range("J5").Value = 5.08
$ws.Range("K5").Value = 2.87
$ws.Range("L5").Value = 6.12
$ws.Range("M5").Value = 4.62
$ws.Range("E6").Value = 6.3
$ws.Range("H6").Value = 6.54
$ws.Range("I6").Value = 6.3
$ws.Range("J6").Value = 5.43
$ws.Range("L6").Value = 6.54
$ws.Range("M6").Value = 4.38
$ws.Range("E7").Value = 6.14
$ws.Range("I7").Value = 6.14
$ws.Range("K7").Value = 2.28
$ws.Range("M7").Value = 4.27
$ws.Range("E8").Value = 6.07
$ws.Range("H8").Value = 6.74
$ws.Range("I8").Value = 6.07
$ws.Range("J8").Value = 5.59
$ws.Range("K8").Value = 2.63
$ws.Range("L8").Value = 6.74
$ws.Range("M8").Value = 4.47
$ws.Range("E9").Value = 6.05
$ws.Range("H9").Value = 6.77
$ws.Range("I9").Value = 6.05
$ws.Range("J9").Value = 5.62
$ws.Range("K9").Value = 2.44
$ws.Range("L9").Value = 6.77
$ws.Range("M9").Value = 4.36
$ws.Range("E10").Value = 9.720000000000001
$ws.Range("H10").Value = 4.58
$ws.Range("I10").Value = 9.720000000000001
$ws.Range("J10").Value = 3.8
$ws.Range("K10").Value = 2.88
$ws.Range("L10").Value = 4.58
$ws.Range("M10").Value = 4.62
$ws.Range("E11").Value = 10.47
$ws.Range("I11").Value = 10.47
$ws.Range("K11").Value = 3.57
$ws.Range("M11").Value = 5.08
$ws.Range("E12").Value = 13.19
$ws.Range("H12").Value = 4.18
$ws.Range("I12").Value = 13.19
$ws.Range("J12").Value = 3.47
$ws.Range("K12").Value = 4.02
$ws.Range("L12").Value = 4.18
$ws.Range("M12").Value = 5.41
$ws.Range("E13").Value = 9.76
$ws.Range("H13").Value = 4.77
$ws.Range("I13").Value = 9.76
$ws.Range("J13").Value = 3.96
$ws.Range("K13").Value = 3.42
$ws.Range("L13").Value = 4.77
$ws.Range("M13").Value = 4.98
$ws.Range("E14").Value = 11.78
$ws.Range("H14").Value = 4.14
$ws.Range("I14").Value = 11.78
$ws.Range("J14").Value = 3.44
$ws.Range("K14").Value = 3.18
$ws.Range("L14").Value = 4.14
$ws.Range("M14").Value = 4.81
$ws.Range("E15").Value = 236.78
$ws.Range("H15").Value = 2.81
$ws.Range("I15").Value = 236.78
$ws.Range("J15").Value = 2.33
$ws.Range("K15").Value = 4.95
$ws.Range("L15").Value = 2.81
$ws.Range("M15").Value = 6.13
$ws.Range("E16").Value = 5.81
$ws.Range("H16").Value = 7.17
$ws.Range("I16").Value = 5.81
$ws.Range("J16").Value = 5.95
$ws.Range("K16").Value = 2.32
$ws.Range("L16").Value = 7.17
$ws.Range("M16").Value = 4.29
$ws.Range("E17").Value = 6.07
$ws.Range("H17").Value = 6.96
$ws.Range("I17").Value = 6.07
$ws.Range("J17").Value = 5.78
$ws.Range("K17").Value = 2.61
$ws.Range("L17").Value = 6.96
$ws.Range("M17").Value = 4.46
$ws.Range("E18").Value = 5.91
$ws.Range("H18").Value = 7.05
$ws.Range("I18").Value = 5.91
$ws.Range("J18").Value = 5.85
$ws.Range("K18").Value = 2.42
$ws.Range("L18").Value = 7.05
$ws.Range("M18").Value = 4.35
$ws.Range("E19").Value = 5.7
$ws.Range("H19").Value = 7.26
$ws.Range("I19").Value = 5.7
$ws.Range("J19").Value = 6.02
$ws.Range("K19").Value = 2.36
$ws.Range("L19").Value = 7.26
$ws.Range("M19").Value = 4.32
$ws.Range("E20").Value = 5.76
$ws.Range("H20").Value = 7.15
$ws.Range("I20").Value = 5.76
$ws.Range("J20").Value = 5.93
$ws.Range("K20").Value = 2.23
$ws.Range("L20").Value = 7.15
$ws.Range("M20").Value = 4.24
$ws.Range("E21").Value = 5.64
$ws.Range("I21").Value = 5.64
$ws.Range("J21").Value = 6.17
$ws.Range("K21").Value = 2.26
$ws.Range("M21").Value = 4.26
